$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells in row 1, columns S and T
$ws.Range("S1").Value = "Time training"
$ws.Range("T1").Value = "Time test"

# Copy the style of the existing header cells (R1) onto the new ones
$ws.Range("R1").Copy()
$ws.Range("S1:T1").PasteSpecial(-4122) # xlPasteFormats

# Update the selected cell to T2
$ws.Range("T2").Select()
